$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 208
$ws.Range("I33").Value = 208
$ws.Range("K33").Value = 208
$ws.Range("M33").Value = 21
$ws.Range("H106").Value = 10860.5
$ws.Range("I106").Value = 3669.5454
$ws.Range("J106").Value = 22160.572
$ws.Range("K106").Value = 3669.5454
$ws.Range("L106").Value = 22160.572
$ws.Range("M106").Value = -3038.5454
$ws.Range("N106").Value = -23422.572
$ws.Range("H116").Value = 15929.429
$ws.Range("I116").Value = 13333.333
$ws.Range("K116").Value = 13333.333
$ws.Range("M116").Value = -9891.333000000001
$ws.Range("H129").Value = 2399.6875
$ws.Range("I129").Value = 2049.75
$ws.Range("J129").Value = 2749.625
$ws.Range("K129").Value = 6149.25
$ws.Range("L129").Value = 8248.875
$ws.Range("M129").Value = -1149.25
$ws.Range("N129").Value = -18248.875
$ws.Range("H132").Value = 3113.5789
$ws.Range("I132").Value = 2884.5625
$ws.Range("K132").Value = 8653.6875
$ws.Range("M132").Value = -6123.6875
$ws.Range("H137").Value = 2685.4
$ws.Range("I137").Value = 2124.2
$ws.Range("J137").Value = 3246.6
$ws.Range("K137").Value = 6372.599999999999
$ws.Range("L137").Value = 9739.799999999999
$ws.Range("M137").Value = -3822.599999999999
$ws.Range("N137").Value = -14839.8
$ws.Range("H140").Value = 69997.71000000001
$ws.Range("J140").Value = 69997.71000000001
$ws.Range("L140").Value = 69997.71000000001
$ws.Range("N140").Value = -80357.71000000001
$ws.Range("H141").Value = 5799.2144
$ws.Range("I141").Value = 3190.6155
$ws.Range("K141").Value = 9571.8465
$ws.Range("M141").Value = -4391.8465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 11685.579
$ws.Range("I2").Value = 814.0833
$ws.Range("J2").Value = 30322.428
$ws.Range("K2").Value = 814.0833
$ws.Range("L2").Value = 30322.428
$ws.Range("M2").Value = -701.0833
$ws.Range("N2").Value = -30548.428
$ws.Range("H32").Value = 3291.2104
$ws.Range("I32").Value = 2236.375
$ws.Range("K32").Value = 2236.375
$ws.Range("M32").Value = -1949.375
$ws.Range("H45").Value = 47622292
$ws.Range("I45").Value = 100001070
$ws.Range("K45").Value = 100001070
$ws.Range("M45").Value = -100000693
$ws.Range("H61").Value = 2781.121
$ws.Range("I61").Value = 2266.6072
$ws.Range("K61").Value = 2266.6072
$ws.Range("M61").Value = -2054.6072
$ws.Range("H80").Value = 54500
$ws.Range("I80").Value = 54500
$ws.Range("K80").Value = 54500
$ws.Range("M80").Value = -53502
$ws.Range("H83").Value = 54500
$ws.Range("I83").Value = 54500
$ws.Range("K83").Value = 163500
$ws.Range("M83").Value = -158508
$ws.Range("H86").Value = 51000
$ws.Range("I86").Value = 51000
$ws.Range("K86").Value = 51000
$ws.Range("M86").Value = -49814
$ws.Range("H88").Value = 2264.6667
$ws.Range("I88").Value = 2110.2222
$ws.Range("J88").Value = 2496.3333
$ws.Range("K88").Value = 2110.2222
$ws.Range("L88").Value = 2496.3333
$ws.Range("M88").Value = -1704.2222
$ws.Range("N88").Value = -3308.3333
$ws.Range("H89").Value = 51000
$ws.Range("I89").Value = 51000
$ws.Range("K89").Value = 153000
$ws.Range("M89").Value = -147072
$ws.Range("H91").Value = 2264.6667
$ws.Range("I91").Value = 2110.2222
$ws.Range("J91").Value = 2496.3333
$ws.Range("K91").Value = 2110.2222
$ws.Range("L91").Value = 2496.3333
$ws.Range("M91").Value = -706.2222000000002
$ws.Range("N91").Value = -5304.3333
$ws.Range("H92").Value = 37000
$ws.Range("J92").Value = 37000
$ws.Range("L92").Value = 37000
$ws.Range("N92").Value = -41992
$ws.Range("H116").Value = 11685.579
$ws.Range("I116").Value = 814.0833
$ws.Range("J116").Value = 30322.428
$ws.Range("K116").Value = 814.0833
$ws.Range("L116").Value = 30322.428
$ws.Range("M116").Value = 1479.9167
$ws.Range("N116").Value = -34910.428
$ws.Range("H132").Value = 7070.091
$ws.Range("I132").Value = 3993.8333
$ws.Range("K132").Value = 11981.4999
$ws.Range("M132").Value = -9451.499899999999
$ws.Range("H136").Value = 2781.121
$ws.Range("I136").Value = 2266.6072
$ws.Range("K136").Value = 6799.821599999999
$ws.Range("M136").Value = -4249.821599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 11685.579
$ws.Range("I3").Value = 814.0833
$ws.Range("J3").Value = 30322.428
$ws.Range("K3").Value = 814.0833
$ws.Range("L3").Value = 30322.428
$ws.Range("M3").Value = -700.0833
$ws.Range("N3").Value = -30550.428
$ws.Range("H54").Value = 1925
$ws.Range("I54").Value = 1925
$ws.Range("K54").Value = 1925
$ws.Range("M54").Value = -1441
$ws.Range("H64").Value = 2177.2222
$ws.Range("I64").Value = 1682.3334
$ws.Range("J64").Value = 2424.6667
$ws.Range("K64").Value = 1682.3334
$ws.Range("L64").Value = 2424.6667
$ws.Range("M64").Value = -1457.3334
$ws.Range("N64").Value = -2874.6667
$ws.Range("H67").Value = 2177.2222
$ws.Range("I67").Value = 1682.3334
$ws.Range("J67").Value = 2424.6667
$ws.Range("K67").Value = 1682.3334
$ws.Range("L67").Value = 2424.6667
$ws.Range("M67").Value = -902.3334
$ws.Range("N67").Value = -3984.6667
$ws.Range("H86").Value = 4593.857
$ws.Range("I86").Value = 2191.6667
$ws.Range("J86").Value = 19007
$ws.Range("K86").Value = 2191.6667
$ws.Range("L86").Value = 19007
$ws.Range("M86").Value = -1068.6667
$ws.Range("N86").Value = -21253
$ws.Range("H89").Value = 4593.857
$ws.Range("I89").Value = 2191.6667
$ws.Range("J89").Value = 19007
$ws.Range("K89").Value = 10958.3335
$ws.Range("L89").Value = 95035
$ws.Range("M89").Value = -5342.333500000001
$ws.Range("N89").Value = -106267
$ws.Range("H99").Value = 2307.1428
$ws.Range("I99").Value = 1733.3334
$ws.Range("J99").Value = 2737.5
$ws.Range("K99").Value = 1733.3334
$ws.Range("L99").Value = 2737.5
$ws.Range("M99").Value = -235.3334
$ws.Range("N99").Value = -5733.5
$ws.Range("H134").Value = 3791.9583
$ws.Range("I134").Value = 2189.5715
$ws.Range("K134").Value = 6568.7145
$ws.Range("M134").Value = -4033.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3514.8823
$ws.Range("I16").Value = 3794.6667
$ws.Range("K16").Value = 3794.6667
$ws.Range("M16").Value = -3507.6667
$ws.Range("H113").Value = 3514.8823
$ws.Range("I113").Value = 3794.6667
$ws.Range("K113").Value = 3794.6667
$ws.Range("M113").Value = -1624.6667
$ws.Range("H132").Value = 5007.636
$ws.Range("I132").Value = 4953.143
$ws.Range("K132").Value = 14859.429
$ws.Range("M132").Value = -12329.429
$ws.Range("H134").Value = 2873.5806
$ws.Range("I134").Value = 2118.76
$ws.Range("J134").Value = 6018.6665
$ws.Range("K134").Value = 6356.280000000001
$ws.Range("L134").Value = 18055.9995
$ws.Range("M134").Value = -3821.280000000001
$ws.Range("N134").Value = -23125.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3191
$ws.Range("J39").Value = 3528.4285
$ws.Range("L39").Value = 10585.2855
$ws.Range("N39").Value = -11173.2855
$ws.Range("H114").Value = 1053.7
$ws.Range("I114").Value = 739.7143
$ws.Range("J114").Value = 1786.3334
$ws.Range("K114").Value = 2219.1429
$ws.Range("L114").Value = 5359.0002
$ws.Range("M114").Value = 1034.8571
$ws.Range("N114").Value = -11867.0002
$ws.Range("H137").Value = 4133.1816
$ws.Range("I137").Value = 1925
$ws.Range("J137").Value = 10021.667
$ws.Range("K137").Value = 5775
$ws.Range("L137").Value = 30065.001
$ws.Range("M137").Value = -675
$ws.Range("N137").Value = -40265.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1170.6316
$ws.Range("I107").Value = 562.5
$ws.Range("K107").Value = 562.5
$ws.Range("M107").Value = 1357.5
$ws.Range("H113").Value = 3338.9
$ws.Range("I113").Value = 1981.6666
$ws.Range("K113").Value = 1981.6666
$ws.Range("M113").Value = 188.3334
$ws.Range("H132").Value = 108457.4
$ws.Range("I132").Value = 171760.17
$ws.Range("K132").Value = 515280.51
$ws.Range("M132").Value = -512750.51

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1387.2646
$ws.Range("I55").Value = 751.7646999999999
$ws.Range("J55").Value = 2022.7646
$ws.Range("K55").Value = 751.7646999999999
$ws.Range("L55").Value = 2022.7646
$ws.Range("M55").Value = -578.7646999999999
$ws.Range("N55").Value = -2368.7646
$ws.Range("H100").Value = 5825.3335
$ws.Range("I100").Value = 4000
$ws.Range("K100").Value = 4000
$ws.Range("M100").Value = -3459
$ws.Range("H122").Value = 5879.64
$ws.Range("I122").Value = 5380.2856
$ws.Range("K122").Value = 16140.8568
$ws.Range("M122").Value = -13690.8568
$ws.Range("H132").Value = 4348.6
$ws.Range("I132").Value = 3167.5454
$ws.Range("J132").Value = 5276.5713
$ws.Range("K132").Value = 9502.636200000001
$ws.Range("L132").Value = 15829.7139
$ws.Range("M132").Value = -6972.636200000001
$ws.Range("N132").Value = -20889.7139
$ws.Range("H136").Value = 4108.2383
$ws.Range("I136").Value = 2137.9
$ws.Range("K136").Value = 6413.700000000001
$ws.Range("M136").Value = -3863.700000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 42500
$ws.Range("J95").Value = 42500
$ws.Range("L95").Value = 42500
$ws.Range("N95").Value = -47992
$ws.Range("H126").Value = 2610.3333
$ws.Range("I126").Value = 1332.1666
$ws.Range("K126").Value = 3996.4998
$ws.Range("M126").Value = -1526.4998
$ws.Range("H132").Value = 3660.9143
$ws.Range("I132").Value = 3776.7932
$ws.Range("K132").Value = 11330.3796
$ws.Range("M132").Value = -8800.3796
$ws.Range("H136").Value = 3334.7896
$ws.Range("J136").Value = 19002.5
$ws.Range("L136").Value = 57007.5
$ws.Range("N136").Value = -62107.5
